# Applies the changes described by the target diff:
#   1. The one table on the deck whose style is the custom "Table_0" style
#      ({BA12230F-CBC3-4B2B-954A-A96C8C59C4AB}) is switched to the built-in
#      table style {47384422-E73B-433A-AF59-FB0D7B21ECC8}.
#   2. The presentation's design/theme palette is changed from the
#      "Integral" colour scheme to the standard "Office Theme" colour
#      scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table -------------------------------------------------
$oldStyleId = "{BA12230F-CBC3-4B2B-954A-A96C8C59C4AB}"
$newStyleId = "{47384422-E73B-433A-AF59-FB0D7B21ECC8}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            if ($shape.Table.Style -eq $oldStyleId) {
                $shape.Table.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Swap the theme colour scheme over to "Office Theme" ---------------
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $hex = $officeColors[$i]
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $themeColors.Colors($i + 1).RGB = $rgb
}

Write-Output "Applied table style + theme colour updates"
